$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear B9's value ("Swag Labs"), move it up to A9 (replacing "Title="),
# and give A9 left-aligned + wrap-text formatting.
$ws.Range("A9").Value2 = $ws.Range("B9").Value2
$ws.Range("B9").Value2 = $null

$ws.Range("A9").HorizontalAlignment = -4131
$ws.Range("A9").WrapText = $true

# Move the active selection from B9 to A9.
$ws.Range("A9").Select()
